# CIERRE 31 DIC 2021
# Roll the payroll sheet from "SEMANA 52" (20-26 Dic 2021) to
# "SEMANA 53" (27 Dic 2021 - 02 Ene 202..) and update the week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header / period label (drives H9, B27, H27, B43, H43, B60 via formulas)
$ws.Range("B9").Value = "SEMANA   53  DEL    27      Al   02   DE   E N E R O          202"

# Keep the active selection where the author left it
$ws.Range("H6").Select() | Out-Null

# --- First block (rows 3-7) ---
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 2200
$ws.Range("K4").Value = 0

# --- Second block (rows 20-26) ---
$ws.Range("K21").Value = 840

# --- Third block (rows 36-41) ---
$ws.Range("K39").Value = 0
